$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 349, pushing the existing rows 349-467 down to 350-468.
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new weekly record.
$ws.Range("A349").Value = 6
$ws.Range("B349").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 45215
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112026
$ws.Range("G349").Value = "Haba"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 480
$ws.Range("K349").Value = 10000
$ws.Range("L349").Value = 12000
$ws.Range("M349").Value = 10542
$ws.Range("N349").Value = "$/saco 25 kilos"
$ws.Range("O349").Value = "Región Metropolitana"
$ws.Range("P349").Value = 422
$ws.Range("Q349").Value = 25
$ws.Range("R349").Value = "Hortaliza"
